# Auto-generated edit script: apply updated market-board price/profit figures
# across the Ravana_Profits workbook sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1883.7
$ws.Range("I137").Value = 1262.6428
$ws.Range("K137").Value = 3787.9284
$ws.Range("M137").Value = -1237.9284
$ws.Range("H138").Value = 7468.4165
$ws.Range("J138").Value = 7988
$ws.Range("L138").Value = 23964
$ws.Range("N138").Value = -34244

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3267.3225
$ws.Range("I32").Value = 2709.6
$ws.Range("K32").Value = 2709.6
$ws.Range("M32").Value = -2422.6

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1077.7142
$ws.Range("I94").Value = 836
$ws.Range("K94").Value = 836
$ws.Range("M94").Value = -385

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4891.4165
$ws.Range("I132").Value = 4283.6665
$ws.Range("K132").Value = 12850.9995
$ws.Range("M132").Value = -10320.9995
$ws.Range("H134").Value = 13547.4
$ws.Range("I134").Value = 13547.4
$ws.Range("K134").Value = 40642.2
$ws.Range("M134").Value = -38107.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.666668
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = $null
$ws.Range("H5").Value = 866.1875
$ws.Range("I5").Value = 696.3333
$ws.Range("K5").Value = 2088.9999
$ws.Range("M5").Value = -1976.9999
$ws.Range("H21").Value = 3000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = -2827
$ws.Range("N21").Value = -15346
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 6000
$ws.Range("M31").Value = -5712
$ws.Range("H33").Value = 1412.5
$ws.Range("I33").Value = 224
$ws.Range("J33").Value = 2601
$ws.Range("K33").Value = 1344
$ws.Range("L33").Value = 15606
$ws.Range("M33").Value = -1061
$ws.Range("N33").Value = -16172
$ws.Range("H108").Value = 827
$ws.Range("I108").Value = 827
$ws.Range("K108").Value = 2481
$ws.Range("M108").Value = 399
$ws.Range("H109").Value = 965.3333
$ws.Range("I109").Value = 900
$ws.Range("K109").Value = 2700
$ws.Range("M109").Value = -1660
$ws.Range("H111").Value = 822
$ws.Range("I111").Value = 822
$ws.Range("K111").Value = 2466
$ws.Range("M111").Value = 601
$ws.Range("H112").Value = 4000
$ws.Range("I112").Value = 4000
$ws.Range("K112").Value = 12000
$ws.Range("M112").Value = -10892
$ws.Range("H113").Value = 740.7857
$ws.Range("J113").Value = 839.55554
$ws.Range("L113").Value = 2518.66662
$ws.Range("N113").Value = -6858.66662
$ws.Range("H114").Value = 2943.875
$ws.Range("I114").Value = 3098
$ws.Range("K114").Value = 9294
$ws.Range("M114").Value = -6040
$ws.Range("H116").Value = 4843.391
$ws.Range("I116").Value = 3449
$ws.Range("K116").Value = 10347
$ws.Range("M116").Value = -6905
$ws.Range("H117").Value = 1549.6
$ws.Range("I117").Value = 1687.5
$ws.Range("K117").Value = 5062.5
$ws.Range("M117").Value = -1620.5
$ws.Range("H118").Value = 4964.3335
$ws.Range("I118").Value = 4000
$ws.Range("J118").Value = 4997.5864
$ws.Range("K118").Value = 12000
$ws.Range("L118").Value = 14992.7592
$ws.Range("M118").Value = -10757
$ws.Range("N118").Value = -17478.7592
$ws.Range("H119").Value = 4000
$ws.Range("I119").Value = 4000
$ws.Range("K119").Value = 12000
$ws.Range("M119").Value = -7162
$ws.Range("H121").Value = 995.13336
$ws.Range("I121").Value = 896.3333
$ws.Range("J121").Value = 1061
$ws.Range("K121").Value = 2688.9999
$ws.Range("L121").Value = 3183
$ws.Range("M121").Value = -1378.9999
$ws.Range("N121").Value = -5803
$ws.Range("H122").Value = 3580.9583
$ws.Range("I122").Value = 340
$ws.Range("J122").Value = 3649.9148
$ws.Range("K122").Value = 3060
$ws.Range("L122").Value = 32849.2332
$ws.Range("M122").Value = -610
$ws.Range("N122").Value = -37749.2332
$ws.Range("H123").Value = 5837.25
$ws.Range("I123").Value = 2745
$ws.Range("J123").Value = 6000
$ws.Range("K123").Value = 8235
$ws.Range("L123").Value = 18000
$ws.Range("M123").Value = -5785
$ws.Range("N123").Value = -22900
$ws.Range("H124").Value = 4893
$ws.Range("I124").Value = 3930
$ws.Range("K124").Value = 11790
$ws.Range("M124").Value = -6880
$ws.Range("H125").Value = 7352.9414
$ws.Range("I125").Value = 5000
$ws.Range("J125").Value = 7500
$ws.Range("K125").Value = 15000
$ws.Range("L125").Value = 22500
$ws.Range("M125").Value = -10080
$ws.Range("N125").Value = -32340
$ws.Range("H126").Value = 6611.706
$ws.Range("I126").Value = 2466.3333
$ws.Range("K126").Value = 7398.999899999999
$ws.Range("M126").Value = -2458.999899999999
$ws.Range("H129").Value = 499.5
$ws.Range("I129").Value = 499.5
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1498.5
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = 3501.5
$ws.Range("H130").Value = 9996.666999999999
$ws.Range("I130").Value = 9990
$ws.Range("K130").Value = 29970
$ws.Range("M130").Value = -24950
$ws.Range("H131").Value = 499.5
$ws.Range("I131").Value = 499.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1498.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = 3541.5
$ws.Range("H132").Value = 4149.5
$ws.Range("I132").Value = 2359.6
$ws.Range("J132").Value = 5428
$ws.Range("K132").Value = 21236.4
$ws.Range("L132").Value = 48852
$ws.Range("M132").Value = -18706.4
$ws.Range("N132").Value = -53912
$ws.Range("H134").Value = 1265
$ws.Range("I134").Value = 1265
$ws.Range("K134").Value = 3795
$ws.Range("M134").Value = 1275
$ws.Range("H135").Value = 866.1875
$ws.Range("I135").Value = 696.3333
$ws.Range("K135").Value = 6266.9997
$ws.Range("M135").Value = -3731.9997
$ws.Range("H137").Value = 8499.75
$ws.Range("J137").Value = 9666.333000000001
$ws.Range("L137").Value = 28998.999
$ws.Range("N137").Value = -39198.999
$ws.Range("H138").Value = 2000
$ws.Range("I138").Value = 2000
$ws.Range("K138").Value = 6000
$ws.Range("M138").Value = -860
$ws.Range("H139").Value = 498
$ws.Range("I139").Value = 498
$ws.Range("K139").Value = 1494
$ws.Range("M139").Value = 3646
$ws.Range("H140").Value = 2632
$ws.Range("I140").Value = 2632
$ws.Range("K140").Value = 7896
$ws.Range("M140").Value = -2716
$ws.Range("H141").Value = 5349
$ws.Range("I141").Value = 5349
$ws.Range("K141").Value = 16047
$ws.Range("M141").Value = -10867

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1424
$ws.Range("I113").Value = 1432
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1432
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 738
$ws.Range("N113").Value = -5740
$ws.Range("H122").Value = 799.75
$ws.Range("I122").Value = 733.3333
$ws.Range("K122").Value = 2199.9999
$ws.Range("M122").Value = 250.0001000000002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2234.0789
$ws.Range("I132").Value = 1907.1
$ws.Range("K132").Value = 5721.299999999999
$ws.Range("M132").Value = -3191.299999999999
